# Added UpdateItemPrice Column in Import / export sheet for Catalog page. [TF00287]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1 (shared string "UpdateItemPrice"), same style as before (s="10")
$ws.Range("G1").Value = "UpdateItemPrice"

# Author's last selection moved to the newly added header cell
$null = $ws.Range("G1").Select()
